$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two new columns before column E ("Client"), pushing the
#    Client..Tier block from E:M out to G:O.
$ws.Range("E1:F1").EntireColumn.Insert()

# 2. Header text for the two freshly-inserted columns.
$ws.Range("E1").Value = "Typist"
$ws.Range("F1").Value = "Typist QC"

# 3. Column width tweaks that accompanied the edit.
$ws.Columns("C").ColumnWidth = 20.56
$ws.Columns("J").ColumnWidth = 19.11

# 4. Add the new order as row 3. Clone row 2's per-column formatting first
#    (same style id in every column) so the new row visually matches, then
#    overwrite the values with the new order's data.
$ws.Range("A2:O2").Copy($ws.Range("A3"))

$ws.Range("A3").Value = 45439.083333333336
$ws.Range("B3").Value = "Be18-002"
$ws.Range("C3").ClearContents()
$ws.Range("D3").ClearContents()
$ws.Range("E3").Value = "SIPL0102"
$ws.Range("F3").Value = "SIPL5317"
$ws.Range("G3").Value = "Beeline"
$ws.Range("H3").Value = "Title"
$ws.Range("I3").Value = "Typing"
$ws.Range("J3").Value = "Commitment Typing"
$ws.Range("K3").Value = "FL"
$ws.Range("L3").Value = "Clay"
$ws.Range("M3").Value = "FLClay"
$ws.Range("N3").Value = "Typing"
$ws.Range("O3").Value = "Typing(T1)"

# E3's box loses its top edge (new border/style combination).
$ws.Range("E3").Borders.Item(8).LineStyle = -4142

$ws.Range("J6").Select()
